# Rename the three worksheets (remove Vietnamese diacritics / spaces -> snake_case)
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "gia_bot_sat"
$wb.Worksheets.Item(2).Name = "gia_than_coc"
$wb.Worksheets.Item(3).Name = "gia_cot_thep"

# "Giá Cốt Thép" (now "gia_cot_thep") stays the active/selected sheet, but the
# in-sheet selection moved from D1 to I22.
$ws3 = $wb.Worksheets.Item("gia_cot_thep")
$ws3.Activate()
$ws3.Range("I22").Select()
